$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 14-17 (previously present combinations with Target cluster = Resolving-Mac for senders,
# and all rows for Sending cluster = Resolving-Mac -> rows covering old A14:T17 block)
$ws.Range("A14:T17").EntireRow.Delete()

# Update remaining data rows (2-13) with recomputed TPM-based values

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Adm"
$ws.Cells.Item(2, 3).Value = "Ramp2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 12.31940633333333
$ws.Cells.Item(2, 8).Value = 36.958219
$ws.Cells.Item(2, 9).Value = 0.3801768364207869
$ws.Cells.Item(2, 10).Value = 0.3801768364207869
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 111.7929483333333
$ws.Cells.Item(2, 14).Value = 335.378845
$ws.Cells.Item(2, 15).Value = 0.7150986351558443
$ws.Cells.Item(2, 16).Value = 0.7150986351558442
$ws.Cells.Item(2, 17).Value = 1377.222755719672
$ws.Cells.Item(2, 18).Value = 12395.00480147705
$ws.Cells.Item(2, 19).Value = 0.2718639368423714
$ws.Cells.Item(2, 20).Value = 0.2718639368423713

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Adm"
$ws.Cells.Item(3, 3).Value = "Ramp2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 12.31940633333333
$ws.Cells.Item(3, 8).Value = 36.958219
$ws.Cells.Item(3, 9).Value = 0.3801768364207869
$ws.Cells.Item(3, 10).Value = 0.3801768364207869
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 31.50896233333333
$ws.Cells.Item(3, 14).Value = 94.52688699999999
$ws.Cells.Item(3, 15).Value = 0.2015513169270731
$ws.Cells.Item(3, 16).Value = 0.2015513169270731
$ws.Cells.Item(3, 17).Value = 388.1717101260281
$ws.Cells.Item(3, 18).Value = 3493.545391134252
$ws.Cells.Item(3, 19).Value = 0.07662514204577804
$ws.Cells.Item(3, 20).Value = 0.07662514204577803

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Adm"
$ws.Cells.Item(4, 3).Value = "Ramp2"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 12.31940633333333
$ws.Cells.Item(4, 8).Value = 36.958219
$ws.Cells.Item(4, 9).Value = 0.3801768364207869
$ws.Cells.Item(4, 10).Value = 0.3801768364207869
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 13.030297
$ws.Cells.Item(4, 14).Value = 39.090891
$ws.Cells.Item(4, 15).Value = 0.08335004791708278
$ws.Cells.Item(4, 16).Value = 0.08335004791708277
$ws.Cells.Item(4, 17).Value = 160.5255233870143
$ws.Cells.Item(4, 18).Value = 1444.729710483129
$ws.Cells.Item(4, 19).Value = 0.03168775753263753
$ws.Cells.Item(4, 20).Value = 0.03168775753263752

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Adm"
$ws.Cells.Item(5, 3).Value = "Ramp2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 18.11265066666666
$ws.Cells.Item(5, 8).Value = 54.33795199999999
$ws.Cells.Item(5, 9).Value = 0.5589563363143816
$ws.Cells.Item(5, 10).Value = 0.5589563363143816
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 111.7929483333333
$ws.Cells.Item(5, 14).Value = 335.378845
$ws.Cells.Item(5, 15).Value = 0.7150986351558443
$ws.Cells.Item(5, 16).Value = 0.7150986351558442
$ws.Cells.Item(5, 17).Value = 2024.866620158381
$ws.Cells.Item(5, 18).Value = 18223.79958142544
$ws.Cells.Item(5, 19).Value = 0.3997089132101254
$ws.Cells.Item(5, 20).Value = 0.3997089132101252

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Adm"
$ws.Cells.Item(6, 3).Value = "Ramp2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 18.11265066666666
$ws.Cells.Item(6, 8).Value = 54.33795199999999
$ws.Cells.Item(6, 9).Value = 0.5589563363143816
$ws.Cells.Item(6, 10).Value = 0.5589563363143816
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 31.50896233333333
$ws.Cells.Item(6, 14).Value = 94.52688699999999
$ws.Cells.Item(6, 15).Value = 0.2015513169270731
$ws.Cells.Item(6, 16).Value = 0.2015513169270731
$ws.Cells.Item(6, 17).Value = 570.7108276128247
$ws.Cells.Item(6, 18).Value = 5136.397448515423
$ws.Cells.Item(6, 19).Value = 0.1126583856888956
$ws.Cells.Item(6, 20).Value = 0.1126583856888956

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Adm"
$ws.Cells.Item(7, 3).Value = "Ramp2"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 18.11265066666666
$ws.Cells.Item(7, 8).Value = 54.33795199999999
$ws.Cells.Item(7, 9).Value = 0.5589563363143816
$ws.Cells.Item(7, 10).Value = 0.5589563363143816
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 13.030297
$ws.Cells.Item(7, 14).Value = 39.090891
$ws.Cells.Item(7, 15).Value = 0.08335004791708278
$ws.Cells.Item(7, 16).Value = 0.08335004791708277
$ws.Cells.Item(7, 17).Value = 236.0132176439146
$ws.Cells.Item(7, 18).Value = 2124.118958795232
$ws.Cells.Item(7, 19).Value = 0.04658903741536075
$ws.Cells.Item(7, 20).Value = 0.04658903741536073

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Adm"
$ws.Cells.Item(8, 3).Value = "Ramp2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.603212
$ws.Cells.Item(8, 8).Value = 4.809636
$ws.Cells.Item(8, 9).Value = 0.04947511672073613
$ws.Cells.Item(8, 10).Value = 0.04947511672073613
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 111.7929483333333
$ws.Cells.Item(8, 14).Value = 335.378845
$ws.Cells.Item(8, 15).Value = 0.7150986351558443
$ws.Cells.Item(8, 16).Value = 0.7150986351558442
$ws.Cells.Item(8, 17).Value = 179.22779628338
$ws.Cells.Item(8, 18).Value = 1613.05016655042
$ws.Cells.Item(8, 19).Value = 0.0353795884411745
$ws.Cells.Item(8, 20).Value = 0.03537958844117449

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Adm"
$ws.Cells.Item(9, 3).Value = "Ramp2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.603212
$ws.Cells.Item(9, 8).Value = 4.809636
$ws.Cells.Item(9, 9).Value = 0.04947511672073613
$ws.Cells.Item(9, 10).Value = 0.04947511672073613
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 31.50896233333333
$ws.Cells.Item(9, 14).Value = 94.52688699999999
$ws.Cells.Item(9, 15).Value = 0.2015513169270731
$ws.Cells.Item(9, 16).Value = 0.2015513169270731
$ws.Cells.Item(9, 17).Value = 50.515546520348
$ws.Cells.Item(9, 18).Value = 454.639918683132
$ws.Cells.Item(9, 19).Value = 0.00997177493018502
$ws.Cells.Item(9, 20).Value = 0.00997177493018502

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Adm"
$ws.Cells.Item(10, 3).Value = "Ramp2"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.603212
$ws.Cells.Item(10, 8).Value = 4.809636
$ws.Cells.Item(10, 9).Value = 0.04947511672073613
$ws.Cells.Item(10, 10).Value = 0.04947511672073613
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 13.030297
$ws.Cells.Item(10, 14).Value = 39.090891
$ws.Cells.Item(10, 15).Value = 0.08335004791708278
$ws.Cells.Item(10, 16).Value = 0.08335004791708277
$ws.Cells.Item(10, 17).Value = 20.890328513964
$ws.Cells.Item(10, 18).Value = 188.012956625676
$ws.Cells.Item(10, 19).Value = 0.00412375334937662
$ws.Cells.Item(10, 20).Value = 0.00412375334937662

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Adm"
$ws.Cells.Item(11, 3).Value = "Ramp2"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.3691416666666666
$ws.Cells.Item(11, 8).Value = 1.107425
$ws.Cells.Item(11, 9).Value = 0.01139171054409548
$ws.Cells.Item(11, 10).Value = 0.01139171054409548
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 111.7929483333333
$ws.Cells.Item(11, 14).Value = 335.378845
$ws.Cells.Item(11, 15).Value = 0.7150986351558443
$ws.Cells.Item(11, 16).Value = 0.7150986351558442
$ws.Cells.Item(11, 17).Value = 41.26743526934721
$ws.Cells.Item(11, 18).Value = 371.4069174241249
$ws.Cells.Item(11, 19).Value = 0.008146196662173118
$ws.Cells.Item(11, 20).Value = 0.008146196662173115

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Adm"
$ws.Cells.Item(12, 3).Value = "Ramp2"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.3691416666666666
$ws.Cells.Item(12, 8).Value = 1.107425
$ws.Cells.Item(12, 9).Value = 0.01139171054409548
$ws.Cells.Item(12, 10).Value = 0.01139171054409548
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 31.50896233333333
$ws.Cells.Item(12, 14).Value = 94.52688699999999
$ws.Cells.Item(12, 15).Value = 0.2015513169270731
$ws.Cells.Item(12, 16).Value = 0.2015513169270731
$ws.Cells.Item(12, 17).Value = 11.63127087066389
$ws.Cells.Item(12, 18).Value = 104.681437835975
$ws.Cells.Item(12, 19).Value = 0.002296014262214468
$ws.Cells.Item(12, 20).Value = 0.002296014262214467

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Adm"
$ws.Cells.Item(13, 3).Value = "Ramp2"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.3691416666666666
$ws.Cells.Item(13, 8).Value = 1.107425
$ws.Cells.Item(13, 9).Value = 0.01139171054409548
$ws.Cells.Item(13, 10).Value = 0.01139171054409548
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 13.030297
$ws.Cells.Item(13, 14).Value = 39.090891
$ws.Cells.Item(13, 15).Value = 0.08335004791708278
$ws.Cells.Item(13, 16).Value = 0.08335004791708277
$ws.Cells.Item(13, 17).Value = 4.810025551741666
$ws.Cells.Item(13, 18).Value = 43.290229965675
$ws.Cells.Item(13, 19).Value = 0.0009494996197078954
$ws.Cells.Item(13, 20).Value = 0.0009494996197078951
